# Criação da tela de cadastro de patente
# Adds a new data row (row 10) to the time-tracking sheet, mirroring the
# date/time formatting already used by the rows above it, and moves the
# selection to C10 (just past the newly entered data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: 15/10/2013, 2h05 (2.0833... hours => 0.0868055... of a day)
$ws.Range("A10").Value = 41562
$ws.Range("B10").Value = 0.086805555555555566

# Match the formatting of the existing rows by copying their styles
# (date style from A9, time style from B9) onto the new cells.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection to C10, as in the final workbook state.
$ws.Range("C10").Select()
